$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Penalty method -> steepest descent method: each penalty value is now
# computed as 0.3 times the previous row's value (instead of 0.5 times).
$prev = $ws.Cells.Item(1, 1).Value2

for ($r = 2; $r -le 5; $r++) {
    $val = $prev * 0.3
    $ws.Cells.Item($r, 1).Value = $val
    $prev = $val
}
